$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first question ("Where can I find access review calendar?") used to
# include a second line ("When the next review kick off?"). That second
# line is being dropped, so the cell becomes a single line of text.
$ws.Range("A1").Value = "Where can I find access review calendar?"

# Row 1 no longer needs the taller (wrapped-text) row height, so auto-fit
# it back down to the sheet's default height now that it is single-line.
$ws.Rows.Item(1).AutoFit()

# Column B widens to fit its (unchanged) longest text now that column A
# no longer dominates the layout.
$ws.Columns.Item(2).ColumnWidth = 47.833333333333336

# Move/select A2 as the active cell (was A9 previously).
$ws.Range("A2").Select()
